$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 21:43"

# Row 4
$ws.Range("B4").Value = 3807731
$ws.Range("C4").Value = 37719
$ws.Range("D4").Value = 1758456
$ws.Range("E4").Value = 1906707
$ws.Range("G4").Value = 504
$ws.Range("H4").Value = 142568

# Row 6
$ws.Range("B6").Value = 1076747
$ws.Range("C6").Value = 36290
$ws.Range("D6").Value = 678010
$ws.Range("E6").Value = 371909
$ws.Range("G6").Value = 543
$ws.Range("H6").Value = 26828

# Row 8 -> Sudafrica
$ws.Range("A8").Value = "Sudafrica"
$ws.Range("B8").Value = 350879
$ws.Range("C8").Value = 13285
$ws.Range("D8").Value = 182230
$ws.Range("E8").Value = 163701
$ws.Range("G8").Value = 144
$ws.Range("H8").Value = 4948

# Row 9 -> Peru
$ws.Range("A9").Value = "Peru"
$ws.Range("B9").Value = 345537
$ws.Range("D9").Value = 233982
$ws.Range("E9").Value = 98756
$ws.Range("H9").Value = 12799

# Row 19
$ws.Range("B19").Value = 202561
$ws.Range("C19").Value = 216
$ws.Range("E19").Value = 5899

# Row 24
$ws.Range("B24").Value = 109993
$ws.Range("C24").Value = 324
$ws.Range("D24").Value = 96907
$ws.Range("E24").Value = 4238
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = 8848

# Row 67
$ws.Range("B67").Value = 16068
$ws.Range("C67").Value = 461
$ws.Range("E67").Value = 6982
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 83

# Row 78
$ws.Range("B78").Value = 10551
$ws.Range("C78").Value = 582
$ws.Range("D78").Value = 2902
$ws.Range("E78").Value = 7595
$ws.Range("G78").Value = 7
$ws.Range("H78").Value = 54

# Row 98
$ws.Range("B98").Value = 4485
$ws.Range("C98").Value = 96
$ws.Range("D98").Value = 1348
$ws.Range("E98").Value = 3082
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 55

# Row 120 -> Cabo Verde
$ws.Range("A120").Value = "Cabo Verde"
$ws.Range("B120").Value = 2014
$ws.Range("C120").Value = 75
$ws.Range("D120").Value = 913
$ws.Range("E120").Value = 1080
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 21

# Row 121 -> Eslovaquia
$ws.Range("A121").Value = "Eslovaquia"
$ws.Range("B121").Value = 1976
$ws.Range("C121").Value = 11
$ws.Range("D121").Value = 1523
$ws.Range("E121").Value = 425
$ws.Range("H121").Value = 28

# Row 122 -> Guinea-Bisau
$ws.Range("A122").Value = "Guinea-Bisau"
$ws.Range("B122").Value = 1949
$ws.Range("C122").Value = 22
$ws.Range("D122").Value = 803
$ws.Range("E122").Value = 1120
$ws.Range("H122").Value = 26

# Row 123 -> Eslovenia
$ws.Range("A123").Value = "Eslovenia"
$ws.Range("B123").Value = 1940
$ws.Range("C123").Value = 24
$ws.Range("D123").Value = 1568
$ws.Range("E123").Value = 261
$ws.Range("H123").Value = 111

# Row 131
$ws.Range("B131").Value = 1581
$ws.Range("C131").Value = 5
$ws.Range("D131").Value = 701
$ws.Range("E131").Value = 437
$ws.Range("G131").Value = 3
$ws.Range("H131").Value = 443

# Row 152
$ws.Range("B152").Value = 743
$ws.Range("C152").Value = 2
$ws.Range("D152").Value = 337
$ws.Range("E152").Value = 392

# Row 155 -> Angola
$ws.Range("A155").Value = "Angola"
$ws.Range("B155").Value = 687
$ws.Range("C155").Value = 49
$ws.Range("D155").Value = 210
$ws.Range("E155").Value = 448
$ws.Range("H155").Value = 29

# Row 156 -> Malta
$ws.Range("A156").Value = "Malta"
$ws.Range("B156").Value = 675
$ws.Range("C156").Value = 1
$ws.Range("D156").Value = 662
$ws.Range("E156").Value = 4
$ws.Range("H156").Value = 9

# Row 171
$ws.Range("B171").Value = 262
$ws.Range("C171").Value = 7
$ws.Range("E171").Value = 149
